# Rooli-excel: Tilaajan urakanvalvoja saa kirjoitusoikeuden vain omaan urakkaan.
#
# Column E ("Tilaajan urakanvalvoja") on the "Oikeudet" sheet currently uses
# "R*,W*[...]" (read/write to ALL contracts) in several rows. Restrict the
# write-access to the user's own contract by dropping the trailing "*" on
# the W-flag: "R*,W*" -> "R*,W" and "R*,W*,sido" -> "R*,W,sido".
# (Read access, "R*", stays "all contracts" - unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")
$ws.Activate()

# Rows whose column-E permission string is the plain "R*,W*" -> "R*,W"
$plainRows = @(7, 8, 9, 10, 11, 12, 20, 21, 22, 23, 24, 25, 27, 28, 30, 32, 33, 53)
foreach ($row in $plainRows) {
    $ws.Range("E$row").Value = "R*,W"
}

# Rows whose column-E permission string carries the extra "sido" flag:
# "R*,W*,sido" -> "R*,W,sido"
$sidoRows = @(29, 31)
foreach ($row in $sidoRows) {
    $ws.Range("E$row").Value = "R*,W,sido"
}

# Row 53 had drifted onto the "alternate band" cell style (borders/centered
# text, but flagged with an extra unused fill-apply bit) instead of the
# plain style the rest of the data rows use. Re-assert the same visual
# formatting so it normalizes back onto the common style bucket, matching
# its sibling cells (D53, F53, ...) and every other row in column E.
$e53 = $ws.Range("E53")
$e53.HorizontalAlignment = -4108   # xlCenter
$e53.Borders.LineStyle = 1         # xlContinuous (thin box border, as before)
$e53.Interior.Pattern = -4142      # xlNone (no fill, as before)

# View tweaks that came along with this edit: a lower zoom level and the
# final cursor position/selection left on the sheet.
$win = $excel.ActiveWindow()
$win.Zoom = 125
$ws.Range("E61").Select() | Out-Null
